# Updates computed market/profit figures on several crafting-leve sheets.
# Generated to match the "chore: update Sheets via scheduled runner" commit.

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $addr = "$col$Row"
        $val = $Values[$col]
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}

# --- ARM sheet ---
Set-Cells "ARM" 18 @{ H = 5000; J = 5000; L = 5000; N = -5644 }
Set-Cells "ARM" 45 @{ H = 669.8889; I = 625.5714; J = 825; K = 625.5714; L = 825; M = -248.5714; N = -1579 }
Set-Cells "ARM" 61 @{ H = 491540.75; I = 386272.3; J = 674006.0600000001; K = 386272.3; L = 674006.0600000001; M = -386060.3; N = -674430.0600000001 }
Set-Cells "ARM" 104 @{ H = 34500; J = 34500; L = 34500; N = -41488 }
Set-Cells "ARM" 110 @{ H = 1114.3077; I = 1158.6; J = 966.6667; K = 1158.6; L = 966.6667; M = 886.4000000000001; N = -5056.6667 }
Set-Cells "ARM" 122 @{ H = 2246.913; I = 2139.1177; J = 2552.3333; K = 6417.353099999999; L = 7656.999899999999; M = -3967.353099999999; N = -12556.9999 }
Set-Cells "ARM" 136 @{ H = 491540.75; I = 386272.3; J = 674006.0600000001; K = 1158816.9; L = 2022018.18; M = -1156266.9; N = -2027118.18 }

# --- BSM sheet ---
Set-Cells "BSM" 26 @{ H = 19188.4; I = 19188.4; K = 19188.4; M = -18896.4 }

# --- CRP sheet ---
Set-Cells "CRP" 36 @{ H = 10700; I = 6400; J = 15000; K = 6400; L = 15000; M = -6012; N = -15776 }
Set-Cells "CRP" 40 @{ H = 10700; I = 6400; J = 15000; K = 6400; L = 15000; M = -6240; N = -15320 }
Set-Cells "CRP" 42 @{ H = 10062; I = 0; J = 10062; K = 0; L = 10062; M = $null; N = -11248 }
Set-Cells "CRP" 44 @{ H = 18500; I = 10000; J = 21333.334; K = 10000; L = 21333.334; M = -9558; N = -22217.334 }
Set-Cells "CRP" 55 @{ H = 25000; J = 25000; L = 25000; N = -25630 }
Set-Cells "CRP" 58 @{ H = 4671.4326; I = 6933.9375; K = 6933.9375; M = -6730.9375 }
Set-Cells "CRP" 99 @{ H = 26403.55; I = 31525.516; J = 2257.1428; K = 31525.516; L = 2257.1428; M = -30027.516; N = -5253.1428 }
Set-Cells "CRP" 107 @{ H = 470.2; I = 459.1111; J = 486.83334; K = 459.1111; L = 486.83334; M = 1460.8889; N = -4326.83334 }
Set-Cells "CRP" 126 @{ H = 26403.55; I = 31525.516; J = 2257.1428; K = 94576.548; L = 6771.428400000001; M = -92106.548; N = -11711.4284 }
Set-Cells "CRP" 136 @{ H = 4671.4326; I = 6933.9375; K = 20801.8125; M = -18251.8125 }

# --- CUL sheet ---
Set-Cells "CUL" 132 @{ H = 1639.1; I = 3344.375; J = 1019; K = 30099.375; L = 9171; M = -27569.375; N = -14231 }
Set-Cells "CUL" 140 @{ H = 7060.4546; I = 724.5185; K = 2173.5555; M = 3006.4445 }

# --- GSM sheet ---
Set-Cells "GSM" 25 @{ H = 69004.5; J = 69004.5; L = 69004.5; N = -70062.5 }
Set-Cells "GSM" 102 @{ H = 3036.1052; I = 2680.5454; J = 3525; K = 2680.5454; L = 3525; M = -1058.5454; N = -6769 }
Set-Cells "GSM" 113 @{ H = 100024610; I = 166701010; J = 10006.5; K = 166701010; L = 10006.5; M = -166698840; N = -14346.5 }
Set-Cells "GSM" 126 @{ H = 3102; I = 2625.4285; J = 3472.6667; K = 7876.2855; L = 10418.0001; M = -5406.2855; N = -15358.0001 }
Set-Cells "GSM" 132 @{ H = 2406995.2; I = 3574205.8; J = 3914.7058; K = 10722617.4; L = 11744.1174; M = -10720087.4; N = -16804.1174 }

# --- LTW sheet ---
Set-Cells "LTW" 3 @{ H = 3680.8; I = 3504; J = 3725; K = 3504; L = 3725; M = -3392; N = -3949 }
Set-Cells "LTW" 15 @{ H = 3680.8; I = 3504; J = 3725; K = 3504; L = 3725; M = -3334; N = -4065 }
Set-Cells "LTW" 133 @{ H = 141375; J = 141375; L = 141375; N = -146435 }

# --- WVR sheet ---
Set-Cells "WVR" 122 @{ H = 1396.5416; I = 1232.8235; J = 1794.1428; K = 3698.4705; L = 5382.428400000001; M = -1248.4705; N = -10282.4284 }
